$d = $word.ActiveDocument

# --- 1. Insert the new paragraphs after the SVG hyperlink paragraph -------
# Locate the SVG hyperlink paragraph and expand to its full paragraph range
# so we can compute the insertion point right after it (and before the
# trailing bookmark paragraph).
$rng = $d.Content
$found = $rng.Find.Execute("https://www.w3schools.com/graphics/svg_intro.asp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the SVG hyperlink paragraph"
}
$rng.Expand(4) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)

$newContentXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
    <w:t>ReactStrap make col full width</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
    <w:instrText xml:space="preserve"> HYPERLINK "https://stackoverflow.com/questions/46151515/react-bootstrap-causing-margins-on-left-and-right-side" </w:instrText>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="3"/>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
    <w:t>https://stackoverflow.com/questions/46151515/react-bootstrap-causing-margins-on-left-and-right-side</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en"/>
    </w:rPr>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertPoint.InsertXML($newContentXml)

# --- 2. Mark the Hyperlink character style as a recommended "quick style" -
# (adds <w:qFormat/> to the style definition, matching w:qFormat="1")
$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.QuickStyle = $true
